{"js": "const body = context.document.body;\nconst results = body.search(\"Bonus\", {matchCase: false});\nresults.load(\"items/text,items/font/color\");\nawait context.sync();\nlet out = [];\nfor (const r of results.items) {\n  out.push(JSON.stringify(r.text) + \" color=\" + r.font.color);\n}\nreturn out.join(\"\\n\");\n", "ps1": "$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"Bonus\"\nif ($find.Execute()) {\n    Write-Output \"found\"\n    Write-Output $find.Parent.Font.Color\n}\n"}
